$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.278.65'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.238.18'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '294.43'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.98'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.514'
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.480'
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.57'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('E12').Value = '  +3.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.56'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '2.583.76'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.90'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('D16').Value = '2.235.38'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.738'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '40.216.84'
$ws.Range('D19').Value = '0.0₃0892'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.46'
$ws.Range('E20').Value = '  +6.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.85'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.74'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.15'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.83'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.88'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.30'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '155.51'
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.39'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.97'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0723'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.91'
$ws.Range('E36').Value = '  +6.21%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.94'
$ws.Range('E38').Value = '  -4.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0983'
$ws.Range('E39').Value = '  -3.24%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = '2.136.53'
$ws.Range('E41').Value = '  +5.28%  '
$ws.Range('E42').Value = '  +1.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '18.30'
$ws.Range('E43').Value = '  +11.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.15'
$ws.Range('E44').Value = '  -3.44%  '
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.84'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('E47').Value = '  +4.41%  '
$ws.Range('D48').Value = '2.448.21'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.73'
$ws.Range('E50').Value = '  -3.03%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '89.30'
$ws.Range('E51').Value = '  -0.83%  '
